$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Seed new rows 16-19 by duplicating the formatting of the existing
#    "Shorts" rows (14 = BLK/shade1, 15 = SND/shade2). Column A (Category),
#    D (Fit), E (Sizes) and F (Shade) already end up with the right values
#    after this copy, since the new SKUs share the same Category/Fit/Sizes
#    and BLK/SND shade pattern as the existing Shorts rows.
# ---------------------------------------------------------------------------
$ws.Range("A14:J15").Copy($ws.Range("A16:J17"))
$ws.Range("A14:J15").Copy($ws.Range("A18:J19"))
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. New product data
#    DL01 / DOUBLE LAYER SHORTS -> rows 16 (BLK) & 17 (SND)
#    CYLS / CYCLING SHORTS      -> rows 18 (BLK) & 19 (SND)
# ---------------------------------------------------------------------------
$img16 = "https://vnhnaiduhall.com/cdn/shop/files/1_2f406cfd-a46a-4866-a225-8158f15ec305.jpg?v=1757271038&width=800,https://vnhnaiduhall.com/cdn/shop/files/4_0b3f815a-227b-47d4-844f-d60ebaddc691.jpg?v=1757271038&width=800"
$img17 = "https://vnhnaiduhall.com/cdn/shop/files/1_5ae150d4-cda3-4849-9b3b-a340261b6f30.jpg?v=1757271038&width=800,https://vnhnaiduhall.com/cdn/shop/files/4_4be02acf-d69a-4234-bbe0-eb1bb6dbc7fc.jpg?v=1757271038&width=800"
$img18 = "https://vnhnaiduhall.com/cdn/shop/files/3_021e435b-7045-4df4-a469-ef77b6a219ca.jpg?v=1757271040&width=800,https://vnhnaiduhall.com/cdn/shop/files/4_f3e232c7-1697-4c24-8387-0bf73b56b752.jpg?v=1757271041&width=800"
$img19 = "https://vnhnaiduhall.com/cdn/shop/files/1_be6f8133-5312-41ad-8f0d-b890cdc9c4db.jpg?v=1757271040&width=800,https://vnhnaiduhall.com/cdn/shop/files/4_ecc4cee6-1ccb-4ee8-b3a0-1d1bc5e50011.jpg?v=1757271040&width=800"

$rows = @(
  @{ Row = 16; Style = "DL01"; StyleName = "DOUBLE LAYER SHORTS"; Img = $img16; MRP = "295"; WSP = "197.65"; DSP = 177.89 },
  @{ Row = 17; Style = "DL01"; StyleName = "DOUBLE LAYER SHORTS"; Img = $img17; MRP = "295"; WSP = "197.65"; DSP = 177.89 },
  @{ Row = 18; Style = "CYLS"; StyleName = "CYCLING SHORTS";      Img = $img18; MRP = "275"; WSP = "184.25"; DSP = 165.83 },
  @{ Row = 19; Style = "CYLS"; StyleName = "CYCLING SHORTS";      Img = $img19; MRP = "275"; WSP = "184.25"; DSP = 165.83 }
)

foreach ($r in $rows) {
  $row = $r.Row

  # Style (B) / StyleName (C) - plain text columns, matching the bold (B)
  # and normal (C) formatting already copied from the source rows.
  $ws.Cells.Item($row, 2).Value = $r.Style
  $ws.Cells.Item($row, 3).Value = $r.StyleName

  # ImageURLs (G) - text content; the hyperlink itself is added below.
  $ws.Cells.Item($row, 7).Value = $r.Img

  # MRP (H) and WSP (I) are stored as *text* in the source data (not
  # numbers), matching the rest of the sheet. Force text entry via a
  # temporary "@" number format, then restore the cell's normal
  # appearance (General / Bold) so the visible formatting still matches
  # the copied rows above.
  $hCell = $ws.Cells.Item($row, 8)
  $hCell.NumberFormat = "@"
  $hCell.Value = $r.MRP
  $hCell.Style = "Normal"

  $iCell = $ws.Cells.Item($row, 9)
  $iCell.NumberFormat = "@"
  $iCell.Value = $r.WSP
  $iCell.Style = "Normal"
  $iCell.Font.Bold = $true

  # DSP (J) - numeric value.
  $ws.Cells.Item($row, 10).Value = $r.DSP

  # Row height to match the rest of the data rows.
  $ws.Rows.Item($row).RowHeight = 15.6
}

# ---------------------------------------------------------------------------
# 3. Hyperlinks on the ImageURLs column, pointing at the same text shown in
#    the cell (same pattern as the existing G14/G15 hyperlinks).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G17"), $img17)
$ws.Range("G17").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("G16"), $img16)
$ws.Range("G16").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("G19"), $img19)
$ws.Range("G19").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("G18"), $img18)
$ws.Range("G18").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4. Column widths (best-fit on the columns that now hold the longest
#    content: Category, StyleName, Fit, ImageURLs).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7.5
$ws.Columns.Item(3).ColumnWidth = 23.666666666666668
$ws.Columns.Item(4).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(7).ColumnWidth = 232.66666666666666

# ---------------------------------------------------------------------------
# 5. Final selection, matching where the cursor ends up after the last
#    edit.
# ---------------------------------------------------------------------------
$ws.Range("K20").Select()
